# Fixed the original data and added output data. Added headings, graphs,
# and preliminary work on DML.
#
# This script adds nine more European countries (Czechia, Ireland, Croatia,
# Cyprus, Hungary, Serbia, Switzerland, Norway, Greece) with their annual
# urban-population-% figures (2014-2022) to the "Data" sheet, directly below
# the existing rows, and updates the sheet's active selection to match the
# authored workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$countries = @(
    "Czechia",
    "Ireland",
    "Croatia",
    "Cyprus",
    "Hungary",
    "Serbia",
    "Switzerland",
    "Norway",
    "Greece"
)

$values = @(
    @(73.384, 73.477000000000004, 73.569999999999993, 73.674999999999997, 73.792000000000002, 73.921000000000006, 74.061000000000007, 74.213999999999999, 74.376999999999995),
    @(62.34, 62.537999999999997, 62.737000000000002, 62.947000000000003, 63.17, 63.405000000000001, 63.652999999999999, 63.911999999999999, 64.183000000000007),
    @(55.923000000000002, 56.155000000000001, 56.402999999999999, 56.667000000000002, 56.947000000000003, 57.241999999999997, 57.552999999999997, 57.878, 58.219000000000001),
    @(67.031000000000006, 66.945999999999998, 66.881, 66.835999999999999, 66.81, 66.805000000000007, 66.820999999999998, 66.855999999999995, 66.911000000000001),
    @(70.227000000000004, 70.5, 70.778999999999996, 71.061999999999998, 71.350999999999999, 71.644000000000005, 71.941999999999993, 72.245000000000005, 72.552000000000007),
    @(55.582999999999998, 55.695999999999998, 55.81, 55.942, 56.091999999999999, 56.26, 56.445999999999998, 56.651000000000003, 56.872999999999998),
    @(73.695999999999998, 73.718000000000004, 73.739000000000004, 73.760999999999996, 73.796999999999997, 73.849000000000004, 73.915000000000006, 73.995999999999995, 74.091999999999999),
    @(80.691999999999993, 81.090999999999994, 81.484999999999999, 81.870999999999995, 82.248000000000005, 82.616, 82.974000000000004, 83.322999999999993, 83.664000000000001),
    @(77.701999999999998, 78.046000000000006, 78.387, 78.724000000000004, 79.058000000000007, 79.388000000000005, 79.715000000000003, 80.037999999999997, 80.356999999999999)
)

$startRow = 20
for ($i = 0; $i -lt $countries.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $countries[$i]
    $rowValues = $values[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $rowValues[$c]
    }
}

# Update the sheet's active cell / selection to match the authored state.
$ws.Range("D15").Select()
